# Auto-generated edit script: updates cryptos list values per commit
# "Updated cryptos list on Sat Mar  9 03:24:13 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.258.28'
$ws.Range('E2').Value = '  +1.71%  '
$ws.Range('D3').Value = '3.927.14'
$ws.Range('E3').Value = '  +0.66%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '483.86'
$ws.Range('E5').Value = '  +4.18%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '146.30'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +1.42%  '
$ws.Range('D7').Value = '0.628'
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').Value = '0.997'
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('D9').Value = '0.728'
$ws.Range('E9').Value = '  -1.07%  '
$ws.Range('D10').Value = '0.171'
$ws.Range('E10').Value = '  +3.74%  '
$ws.Range('D11').Value = '0.0000361'
$ws.Range('E11').Value = '  +5.36%  '
$ws.Range('D12').Value = '42.65'
$ws.Range('E12').Value = '  -0.83%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '10.60'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +2.25%  '
$ws.Range('D14').Value = '4.544.70'
$ws.Range('E14').Value = '  +0.57%  '
$ws.Range('D15').Value = '14.85'
$ws.Range('E15').Value = '  -2.28%  '
$ws.Range('D16').Value = '3.913.20'
$ws.Range('E16').Value = '  -0.50%  '
$ws.Range('E17').Value = '  -0.17%  '
$ws.Range('D18').Value = '19.88'
$ws.Range('E18').Value = '  -0.33%  '
$ws.Range('E19').Value = '  -1.64%  '
$ws.Range('D20').Value = '68.343.04'
$ws.Range('E20').Value = '  +1.60%  '
$ws.Range('D21').Value = '447.48'
$ws.Range('E21').Value = '  +3.37%  '
$ws.Range('D22').Value = '14.83'
$ws.Range('E22').Value = '  +1.31%  '
$ws.Range('D23').Value = '3.37'
$ws.Range('E23').Value = '  +1.20%  '
$ws.Range('D24').Value = '88.75'
$ws.Range('E24').Value = '  +0.02%  '
$ws.Range('D25').Value = '11.39'
$ws.Range('E25').Value = '  +12.92%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '3.60'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +2.73%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = '39.01'
$ws.Range('E27').Value = '  +1.02%  '
$ws.Range('B28').Value = 'RenderToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D28').Value = '10.56'
$ws.Range('E28').Value = '  +9.69%  '
$ws.Range('D29').Value = '5.87'
$ws.Range('E29').Value = '  +2.99%  '
$ws.Range('D30').Value = '699.19'
$ws.Range('E30').Value = '  -5.44%  '
$ws.Range('D31').Value = '13.44'
$ws.Range('E31').Value = '  -1.11%  '
$ws.Range('E32').Value = '  -0.76%  '
$ws.Range('D33').Value = '2.87'
$ws.Range('E33').Value = '  +3.13%  '
$ws.Range('D34').Value = '0.0₃0943'
$ws.Range('E34').Value = '  +20.40%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '41.80'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -2.94%  '
$ws.Range('D36').Value = '59.12'
$ws.Range('E36').Value = '  +1.82%  '
$ws.Range('B37').Value = 'NEARProtocol'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D37').Value = '5.68'
$ws.Range('E37').Value = '  +5.91%  '
$ws.Range('B38').Value = 'Kaspa'
$ws.Range('C38').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.150'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -4.76%  '
$ws.Range('D39').Value = '0.999'
$ws.Range('E39').Value = '  +0.02%  '
$ws.Range('D40').Value = '0.0479'
$ws.Range('E40').Value = '  +0.73%  '
$ws.Range('D41').Value = '2.83'
$ws.Range('E41').Value = '  +14.30%  '
$ws.Range('E42').Value = '  -2.83%  '
$ws.Range('D43').Value = '0.364'
$ws.Range('E43').Value = '  +8.87%  '
$ws.Range('D44').Value = '2.97'
$ws.Range('E44').Value = '  +6.21%  '
$ws.Range('E45').Value = '  +1.61%  '
$ws.Range('E46').Value = '  -0.20%  '
$ws.Range('E47').Value = '  +0.57%  '
$ws.Range('D48').Value = '2.15'
$ws.Range('E48').Value = '  -0.94%  '
$ws.Range('D49').Value = '145.71'
$ws.Range('E49').Value = '  +1.56%  '
$ws.Range('E50').Value = '  -0.76%  '
$ws.Range('D51').Value = '2.85'
$ws.Range('E51').Value = '  -1.67%  '
